$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting the existing rows 56-68 down to 57-69.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with a new weekly entry for
# "Comercializadora del Agro de Limarí" / "Poroto granado".
$ws.Range("A56").Value = 2
$ws.Range("B56").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44588
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100112030
$ws.Range("G56").Value = "Poroto granado"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 500
$ws.Range("K56").Value = 20000
$ws.Range("L56").Value = 23000
$ws.Range("M56").Value = 21500
$ws.Range("N56").Value = "$/malla 25 kilos"
$ws.Range("O56").Value = "Provincia de Limarí"
$ws.Range("P56").Value = 860
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"
